# Append 35 new NBA game rows (rows 557-591) to the "NBA Season" sheet,
# covering game dates 2026-01-08 through 2026-01-12 (Excel serials 46030-46034).
# Columns: A=game_date, B=home_team, C=away_team, D=closing_spread,
#          E=home_score, F=away_score, G=spread_result_difference
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A557").Value = 46030
$ws.Range("A557").NumberFormat = "YYYY-MM-DD"
$ws.Range("B557").Value = "Charlotte Hornets"
$ws.Range("C557").Value = "Indiana Pacers"
$ws.Range("D557").Value = -4.5
$ws.Range("E557").Value = 112
$ws.Range("F557").Value = 114
$ws.Range("G557").Value = -6.5

$ws.Range("A558").Value = 46030
$ws.Range("A558").NumberFormat = "YYYY-MM-DD"
$ws.Range("B558").Value = "Minnesota Timberwolves"
$ws.Range("C558").Value = "Cleveland Cavaliers"
$ws.Range("D558").Value = -3.5
$ws.Range("E558").Value = 131
$ws.Range("F558").Value = 122
$ws.Range("G558").Value = 5.5

$ws.Range("A559").Value = 46030
$ws.Range("A559").NumberFormat = "YYYY-MM-DD"
$ws.Range("B559").Value = "Utah Jazz"
$ws.Range("C559").Value = "Dallas Mavericks"
$ws.Range("D559").Value = 5.5
$ws.Range("E559").Value = 116
$ws.Range("F559").Value = 114
$ws.Range("G559").Value = 7.5

$ws.Range("A560").Value = 46031
$ws.Range("A560").NumberFormat = "YYYY-MM-DD"
$ws.Range("B560").Value = "Boston Celtics"
$ws.Range("C560").Value = "Toronto Raptors"
$ws.Range("D560").Value = -9.5
$ws.Range("E560").Value = 125
$ws.Range("F560").Value = 117
$ws.Range("G560").Value = -1.5

$ws.Range("A561").Value = 46031
$ws.Range("A561").NumberFormat = "YYYY-MM-DD"
$ws.Range("B561").Value = "Orlando Magic"
$ws.Range("C561").Value = "Philadelphia 76ers"
$ws.Range("D561").Value = 3.5
$ws.Range("E561").Value = 91
$ws.Range("F561").Value = 103
$ws.Range("G561").Value = -8.5

$ws.Range("A562").Value = 46031
$ws.Range("A562").NumberFormat = "YYYY-MM-DD"
$ws.Range("B562").Value = "Washington Wizards"
$ws.Range("C562").Value = "New Orleans Pelicans"
$ws.Range("D562").Value = 3.5
$ws.Range("E562").Value = 107
$ws.Range("F562").Value = 128
$ws.Range("G562").Value = -17.5

$ws.Range("A563").Value = 46031
$ws.Range("A563").NumberFormat = "YYYY-MM-DD"
$ws.Range("B563").Value = "Brooklyn Nets"
$ws.Range("C563").Value = "Los Angeles Clippers"
$ws.Range("D563").Value = 5.5
$ws.Range("E563").Value = 105
$ws.Range("F563").Value = 121
$ws.Range("G563").Value = -10.5

$ws.Range("A564").Value = 46031
$ws.Range("A564").NumberFormat = "YYYY-MM-DD"
$ws.Range("B564").Value = "Memphis Grizzlies"
$ws.Range("C564").Value = "Oklahoma City Thunder"
$ws.Range("D564").Value = 5.5
$ws.Range("E564").Value = 116
$ws.Range("F564").Value = 117
$ws.Range("G564").Value = 4.5

$ws.Range("A565").Value = 46031
$ws.Range("A565").NumberFormat = "YYYY-MM-DD"
$ws.Range("B565").Value = "Phoenix Suns"
$ws.Range("C565").Value = "New York Knicks"
$ws.Range("D565").Value = 1.5
$ws.Range("E565").Value = 112
$ws.Range("F565").Value = 107
$ws.Range("G565").Value = 6.5

$ws.Range("A566").Value = 46031
$ws.Range("A566").NumberFormat = "YYYY-MM-DD"
$ws.Range("B566").Value = "Denver Nuggets"
$ws.Range("C566").Value = "Atlanta Hawks"
$ws.Range("D566").Value = 3.5
$ws.Range("E566").Value = 87
$ws.Range("F566").Value = 110
$ws.Range("G566").Value = -19.5

$ws.Range("A567").Value = 46031
$ws.Range("A567").NumberFormat = "YYYY-MM-DD"
$ws.Range("B567").Value = "Golden State Warriors"
$ws.Range("C567").Value = "Sacramento Kings"
$ws.Range("D567").Value = -15.5
$ws.Range("E567").Value = 137
$ws.Range("F567").Value = 103
$ws.Range("G567").Value = 18.5

$ws.Range("A568").Value = 46031
$ws.Range("A568").NumberFormat = "YYYY-MM-DD"
$ws.Range("B568").Value = "Portland Trail Blazers"
$ws.Range("C568").Value = "Houston Rockets"
$ws.Range("D568").Value = 6.5
$ws.Range("E568").Value = 111
$ws.Range("F568").Value = 105
$ws.Range("G568").Value = 12.5

$ws.Range("A569").Value = 46031
$ws.Range("A569").NumberFormat = "YYYY-MM-DD"
$ws.Range("B569").Value = "Los Angeles Lakers"
$ws.Range("C569").Value = "Milwaukee Bucks"
$ws.Range("D569").Value = -3.5
$ws.Range("E569").Value = 101
$ws.Range("F569").Value = 105
$ws.Range("G569").Value = -7.5

$ws.Range("A570").Value = 46032
$ws.Range("A570").NumberFormat = "YYYY-MM-DD"
$ws.Range("B570").Value = "Cleveland Cavaliers"
$ws.Range("C570").Value = "Minnesota Timberwolves"
$ws.Range("D570").Value = -2.5
$ws.Range("E570").Value = 146
$ws.Range("F570").Value = 134
$ws.Range("G570").Value = 9.5

$ws.Range("A571").Value = 46032
$ws.Range("A571").NumberFormat = "YYYY-MM-DD"
$ws.Range("B571").Value = "Indiana Pacers"
$ws.Range("C571").Value = "Miami Heat"
$ws.Range("D571").Value = 6.5
$ws.Range("E571").Value = 123
$ws.Range("F571").Value = 99
$ws.Range("G571").Value = 30.5

$ws.Range("A572").Value = 46032
$ws.Range("A572").NumberFormat = "YYYY-MM-DD"
$ws.Range("B572").Value = "Detroit Pistons"
$ws.Range("C572").Value = "Los Angeles Clippers"
$ws.Range("D572").Value = -1.5
$ws.Range("E572").Value = 92
$ws.Range("F572").Value = 98
$ws.Range("G572").Value = -7.5

$ws.Range("A573").Value = 46032
$ws.Range("A573").NumberFormat = "YYYY-MM-DD"
$ws.Range("B573").Value = "Boston Celtics"
$ws.Range("C573").Value = "San Antonio Spurs"
$ws.Range("D573").Value = 1.5
$ws.Range("E573").Value = 95
$ws.Range("F573").Value = 100
$ws.Range("G573").Value = -3.5

$ws.Range("A574").Value = 46032
$ws.Range("A574").NumberFormat = "YYYY-MM-DD"
$ws.Range("B574").Value = "Chicago Bulls"
$ws.Range("C574").Value = "Dallas Mavericks"
$ws.Range("D574").Value = -4.5
$ws.Range("E574").Value = 125
$ws.Range("F574").Value = 107
$ws.Range("G574").Value = 13.5

$ws.Range("A575").Value = 46032
$ws.Range("A575").NumberFormat = "YYYY-MM-DD"
$ws.Range("B575").Value = "Utah Jazz"
$ws.Range("C575").Value = "Charlotte Hornets"
$ws.Range("D575").Value = 6.5
$ws.Range("E575").Value = 95
$ws.Range("F575").Value = 150
$ws.Range("G575").Value = -48.5

$ws.Range("A576").Value = 46033
$ws.Range("A576").NumberFormat = "YYYY-MM-DD"
$ws.Range("B576").Value = "Orlando Magic"
$ws.Range("C576").Value = "New Orleans Pelicans"
$ws.Range("D576").Value = -7.5
$ws.Range("E576").Value = 128
$ws.Range("F576").Value = 118
$ws.Range("G576").Value = 2.5

$ws.Range("A577").Value = 46033
$ws.Range("A577").NumberFormat = "YYYY-MM-DD"
$ws.Range("B577").Value = "Memphis Grizzlies"
$ws.Range("C577").Value = "Brooklyn Nets"
$ws.Range("D577").Value = -7.5
$ws.Range("E577").Value = 103
$ws.Range("F577").Value = 98
$ws.Range("G577").Value = -2.5

$ws.Range("A578").Value = 46033
$ws.Range("A578").NumberFormat = "YYYY-MM-DD"
$ws.Range("B578").Value = "Toronto Raptors"
$ws.Range("C578").Value = "Philadelphia 76ers"
$ws.Range("D578").Value = 1.5
$ws.Range("E578").Value = 116
$ws.Range("F578").Value = 115
$ws.Range("G578").Value = 2.5

$ws.Range("A579").Value = 46033
$ws.Range("A579").NumberFormat = "YYYY-MM-DD"
$ws.Range("B579").Value = "Portland Trail Blazers"
$ws.Range("C579").Value = "New York Knicks"
$ws.Range("D579").Value = 5.5
$ws.Range("E579").Value = 114
$ws.Range("F579").Value = 123
$ws.Range("G579").Value = -3.5

$ws.Range("A580").Value = 46033
$ws.Range("A580").NumberFormat = "YYYY-MM-DD"
$ws.Range("B580").Value = "Oklahoma City Thunder"
$ws.Range("C580").Value = "Miami Heat"
$ws.Range("D580").Value = -14.5
$ws.Range("E580").Value = 124
$ws.Range("F580").Value = 112
$ws.Range("G580").Value = -2.5

$ws.Range("A581").Value = 46033
$ws.Range("A581").NumberFormat = "YYYY-MM-DD"
$ws.Range("B581").Value = "Minnesota Timberwolves"
$ws.Range("C581").Value = "San Antonio Spurs"
$ws.Range("D581").Value = -2.5
$ws.Range("E581").Value = 104
$ws.Range("F581").Value = 103
$ws.Range("G581").Value = -1.5

$ws.Range("A582").Value = 46033
$ws.Range("A582").NumberFormat = "YYYY-MM-DD"
$ws.Range("B582").Value = "Denver Nuggets"
$ws.Range("C582").Value = "Milwaukee Bucks"
$ws.Range("D582").Value = 6.5
$ws.Range("E582").Value = 108
$ws.Range("F582").Value = 104
$ws.Range("G582").Value = 10.5

$ws.Range("A583").Value = 46033
$ws.Range("A583").NumberFormat = "YYYY-MM-DD"
$ws.Range("B583").Value = "Phoenix Suns"
$ws.Range("C583").Value = "Washington Wizards"
$ws.Range("D583").Value = -15.5
$ws.Range("E583").Value = 112
$ws.Range("F583").Value = 93
$ws.Range("G583").Value = 3.5

$ws.Range("A584").Value = 46033
$ws.Range("A584").NumberFormat = "YYYY-MM-DD"
$ws.Range("B584").Value = "Golden State Warriors"
$ws.Range("C584").Value = "Atlanta Hawks"
$ws.Range("D584").Value = -7.5
$ws.Range("E584").Value = 111
$ws.Range("F584").Value = 124
$ws.Range("G584").Value = -20.5

$ws.Range("A585").Value = 46033
$ws.Range("A585").NumberFormat = "YYYY-MM-DD"
$ws.Range("B585").Value = "Sacramento Kings"
$ws.Range("C585").Value = "Houston Rockets"
$ws.Range("D585").Value = 14.5
$ws.Range("E585").Value = 111
$ws.Range("F585").Value = 98
$ws.Range("G585").Value = 27.5

$ws.Range("A586").Value = 46034
$ws.Range("A586").NumberFormat = "YYYY-MM-DD"
$ws.Range("B586").Value = "Cleveland Cavaliers"
$ws.Range("C586").Value = "Utah Jazz"
$ws.Range("D586").Value = -13.5
$ws.Range("E586").Value = 112
$ws.Range("F586").Value = 123
$ws.Range("G586").Value = -24.5

$ws.Range("A587").Value = 46034
$ws.Range("A587").NumberFormat = "YYYY-MM-DD"
$ws.Range("B587").Value = "Indiana Pacers"
$ws.Range("C587").Value = "Boston Celtics"
$ws.Range("D587").Value = 5.5
$ws.Range("E587").Value = 98
$ws.Range("F587").Value = 96
$ws.Range("G587").Value = 7.5

$ws.Range("A588").Value = 46034
$ws.Range("A588").NumberFormat = "YYYY-MM-DD"
$ws.Range("B588").Value = "Toronto Raptors"
$ws.Range("C588").Value = "Philadelphia 76ers"
$ws.Range("D588").Value = 3.5
$ws.Range("E588").Value = 102
$ws.Range("F588").Value = 115
$ws.Range("G588").Value = -9.5

$ws.Range("A589").Value = 46034
$ws.Range("A589").NumberFormat = "YYYY-MM-DD"
$ws.Range("B589").Value = "Dallas Mavericks"
$ws.Range("C589").Value = "Brooklyn Nets"
$ws.Range("D589").Value = -3.5
$ws.Range("E589").Value = 113
$ws.Range("F589").Value = 105
$ws.Range("G589").Value = 4.5

$ws.Range("A590").Value = 46034
$ws.Range("A590").NumberFormat = "YYYY-MM-DD"
$ws.Range("B590").Value = "Sacramento Kings"
$ws.Range("C590").Value = "Los Angeles Lakers"
$ws.Range("D590").Value = 9.5
$ws.Range("E590").Value = 124
$ws.Range("F590").Value = 112
$ws.Range("G590").Value = 21.5

$ws.Range("A591").Value = 46034
$ws.Range("A591").NumberFormat = "YYYY-MM-DD"
$ws.Range("B591").Value = "Los Angeles Clippers"
$ws.Range("C591").Value = "Charlotte Hornets"
$ws.Range("D591").Value = -4.5
$ws.Range("E591").Value = 117
$ws.Range("F591").Value = 109
$ws.Range("G591").Value = 3.5
